{"js": "// Update the date line and the division problems in the practice-sheet\n// table. Each entry maps the exact existing run text to its replacement;\n// we locate every occurrence via Body.search (exact, case-sensitive,\n// non-wildcard) and swap the text in place with Range.insertText(...,\n// \"Replace\") so the existing run formatting (font/size) is preserved.\nconst replacements = [\n  [\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"],\n  [\"159\u00f75=\", \"230\u00f79=\"],\n  [\"464\u00f78=\", \"316\u00f72=\"],\n  [\"687\u00f73=\", \"266\u00f79=\"],\n  [\"573\u00f77=\", \"856\u00f77=\"],\n  [\"334\u00f73=\", \"233\u00f78=\"],\n  [\"646\u00f73=\", \"910\u00f73=\"],\n  [\"878\u00f78=\", \"558\u00f79=\"],\n  [\"823\u00f78=\", \"421\u00f78=\"],\n  [\"666\u00f74=\", \"194\u00f73=\"],\n  [\"366\u00f78=\", \"782\u00f76=\"],\n  [\"828\u00f76=\", \"222\u00f72=\"],\n  [\"209\u00f74=\", \"582\u00f76=\"],\n  [\"971\u00f72=\", \"522\u00f73=\"],\n  [\"161\u00f79=\", \"108\u00f79=\"],\n  [\"982\u00f77=\", \"280\u00f79=\"],\n  [\"589\u00f77=\", \"759\u00f78=\"],\n  [\"158\u00f76=\", \"334\u00f75=\"],\n  [\"705\u00f78=\", \"871\u00f76=\"],\n  [\"490\u00f77=\", \"133\u00f77=\"],\n  [\"707\u00f77=\", \"964\u00f74=\"],\n  [\"141\u00f75=\", \"123\u00f79=\"],\n  [\"923\u00f75=\", \"483\u00f73=\"],\n  [\"686\u00f73=\", \"391\u00f78=\"],\n  [\"676\u00f75=\", \"346\u00f74=\"],\n  [\"627\u00f79=\", \"629\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the division problems in the practice-sheet\n# table. Each entry maps the exact existing text to its replacement; we\n# drive Word's Find/Replace (wdReplaceAll) over the whole document body\n# for each pair so every matching run is updated in place (preserving\n# the run's existing font/size formatting).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"),\n    @(\"159\u00f75=\", \"230\u00f79=\"),\n    @(\"464\u00f78=\", \"316\u00f72=\"),\n    @(\"687\u00f73=\", \"266\u00f79=\"),\n    @(\"573\u00f77=\", \"856\u00f77=\"),\n    @(\"334\u00f73=\", \"233\u00f78=\"),\n    @(\"646\u00f73=\", \"910\u00f73=\"),\n    @(\"878\u00f78=\", \"558\u00f79=\"),\n    @(\"823\u00f78=\", \"421\u00f78=\"),\n    @(\"666\u00f74=\", \"194\u00f73=\"),\n    @(\"366\u00f78=\", \"782\u00f76=\"),\n    @(\"828\u00f76=\", \"222\u00f72=\"),\n    @(\"209\u00f74=\", \"582\u00f76=\"),\n    @(\"971\u00f72=\", \"522\u00f73=\"),\n    @(\"161\u00f79=\", \"108\u00f79=\"),\n    @(\"982\u00f77=\", \"280\u00f79=\"),\n    @(\"589\u00f77=\", \"759\u00f78=\"),\n    @(\"158\u00f76=\", \"334\u00f75=\"),\n    @(\"705\u00f78=\", \"871\u00f76=\"),\n    @(\"490\u00f77=\", \"133\u00f77=\"),\n    @(\"707\u00f77=\", \"964\u00f74=\"),\n    @(\"141\u00f75=\", \"123\u00f79=\"),\n    @(\"923\u00f75=\", \"483\u00f73=\"),\n    @(\"686\u00f73=\", \"391\u00f78=\"),\n    @(\"676\u00f75=\", \"346\u00f74=\"),\n    @(\"627\u00f79=\", \"629\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, \"wdReplaceAll\")\n}\n"}
